$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New song entry added 2018-11-12: "Wolves" by "Selena Gomez, Marshmello"
# Copy the formatting of the row-index cell above (A44) so the new index
# cell (A45) keeps the same centered/bordered style used throughout column A.
$ws.Range("A44").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("A45").Value = 43

$ws.Range("B45").Value = "Wolves"
$ws.Range("C45").Value = "Selena Gomez, Marshmello"

# Force the date column to stay as literal text (matching the rest of the
# sheet's "d-m-yyyy" text strings) instead of being auto-converted to a
# date serial number, then drop the resulting format override so the cell
# keeps the sheet's default (unstyled) look.
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12-11-2018"
$ws.Range("D45").ClearFormats()
